$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the _GoBack bookmark from its current location (end of
#    the "Under the file menu Method, select Load Method. " paragraph).
# ------------------------------------------------------------------
$existing = $d.Bookmarks("_GoBack")
$existing.Delete()

# ------------------------------------------------------------------
# 2) Split the trailing sentence off the "...and fraction. If not, be
#    sure to label these appropriately. " run into two runs: keep the
#    first sentence as-is and replace the second sentence with new
#    text, comma/period separated for clarity. Track revisions while
#    doing the replacement (and accept immediately after) so the new
#    text is written as its own run with matching rPr instead of being
#    silently coalesced back into the preceding run.
# ------------------------------------------------------------------
$d.TrackRevisions = $true

$target = $d.Content
$target.Find.Execute("If not, be sure to label these appropriately. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Text = "If you load in different files from those listed in the method .xml file, you will need to click the row and edit these labels. "

$d.AcceptAllRevisions() | Out-Null
$d.TrackRevisions = $false

# ------------------------------------------------------------------
# 3) Re-insert the _GoBack bookmark immediately after the new sentence
#    (collapsed bookmark, same as the original). Adding a bookmark
#    whose range sits exactly at the end of a paragraph's text (i.e.
#    right before the paragraph mark) can snap to the start of the
#    document, so we temporarily pad with a throwaway character,
#    anchor the bookmark just before it, then remove the padding.
# ------------------------------------------------------------------
$located = $d.Content
$located.Find.Execute("If you load in different files from those listed in the method .xml file, you will need to click the row and edit these labels. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPos = $located.End

$pad = $d.Range($insertPos, $insertPos)
$pad.InsertAfter("X")

$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$padRange = $d.Range($insertPos, $insertPos + 1)
$padRange.Delete()
